# Generate Report for Handoff
#
# The localization CI run produced a new handoff package for
# "690a0505-d98d-4f23-a886-ab61827ab3f9.md": it is no longer "In
# Translation" -- it is now "Ready for handoff" (machine translation,
# priority "mt"), and fresh handoff timestamps were recorded for both
# target locales (zh-cn / de-de). Reflect that across the Overview
# sheet (row for this file) as well as the per-locale zh-cn and de-de
# detail sheets (row for this file).

$wb = $excel.ActiveWorkbook

# A slightly wider column width, matching the widened "Status" columns
# that now need to fit the longer "Ready for handoff" label.
$newColumnWidth = 16.333333

# --- Overview sheet ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F3").Value = "Ready for handoff"   # de-de status
$wsOverview.Range("G3").Value = "2016-08-16 08:13:34" # Latest HO Xliff Generate Date
$wsOverview.Columns("E:F").ColumnWidth = $newColumnWidth

# --- zh-cn detail sheet --------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"      # Status
$wsZhCn.Range("E3").Value = "mt"                     # Priority
$wsZhCn.Range("H3").Value = "2016-08-16 08:13:29"    # Latest Handoff Datetime
$wsZhCn.Columns("C:C").ColumnWidth = $newColumnWidth

# --- de-de detail sheet --------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"      # Status
$wsDeDe.Range("E3").Value = "mt"                     # Priority
$wsDeDe.Range("H3").Value = "2016-08-16 08:13:34"    # Latest Handoff Datetime
$wsDeDe.Columns("C:C").ColumnWidth = $newColumnWidth
